$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Clear old columns F and G (no longer used)
$ws.Range("F1:G4").Clear()

# Data rows: Índice, Distancia, max, min, Tempo
$data = @(
    @(0, 10917.3,             11579, 9364,  0.209555729230245),
    @(1, 11066.46666666667,   11859, 9396,  0.2522561311721802),
    @(2, 10921.36666666667,   11558, 9991,  0.2720347007115682),
    @(3, 11558.53333333333,   12256, 10469, 0.1916634956995646),
    @(4, 10512.93333333333,   11319, 9611,  0.2410261233647664),
    @(5, 11466.9,             12369, 9922,  0.285493524869283),
    @(6, 11013.23333333333,   11835, 10138, 0.3088875850041707),
    @(7, 10743.9,             11548, 10143, 0.2741494258244833),
    @(8, 10653.63333333333,   11415, 9151,  0.2374705235163371),
    @(9, 10609.66666666667,   11597, 9665,  0.2030163526535034)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}
